$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---- Sheet 1 ----
$ws1.Cells.Item(2,1).Value = "Última actualización: 09:28:24"
$ws1.Cells.Item(3,1).Value = "Total filas: 130"
$ws1.Cells.Item(24,1).Value = "04:17:03"
$ws1.Cells.Item(24,3).Value = "215B_EL PATO"
$ws1.Cells.Item(24,4).Value = 78
$ws1.Cells.Item(25,1).Value = "03:42:43"
$ws1.Cells.Item(25,3).Value = "14_ABASTO"
$ws1.Cells.Item(25,4).Value = 113
$ws1.Cells.Item(55,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(57,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(85,1).Value = "07:17:57"
$ws1.Cells.Item(85,3).Value = "17_ROMERO"
$ws1.Cells.Item(85,4).Value = 96
$ws1.Cells.Item(86,1).Value = "07:50:23"
$ws1.Cells.Item(86,3).Value = "10_OLMOS"
$ws1.Cells.Item(86,4).Value = 63
$ws1.Cells.Item(87,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(88,3).Value = "17_ROMERO"
$ws1.Cells.Item(102,1).Value = "08:52:26"
$ws1.Cells.Item(102,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(102,4).Value = 31
$ws1.Cells.Item(103,1).Value = "07:50:23"
$ws1.Cells.Item(103,3).Value = "17_ROMERO"
$ws1.Cells.Item(103,4).Value = 93
$ws1.Cells.Item(107,1).Value = "09:28:24"
$ws1.Cells.Item(107,4).Value = 5
$ws1.Cells.Item(109,1).Value = "09:28:24"
$ws1.Cells.Item(109,2).Value = "09:35"
$ws1.Cells.Item(109,4).Value = 7
$ws1.Cells.Item(110,2).Value = "09:39"
$ws1.Cells.Item(110,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(110,4).Value = 60
$ws1.Cells.Item(111,2).Value = "09:41"
$ws1.Cells.Item(111,3).Value = "215C_EL PATO"
$ws1.Cells.Item(111,4).Value = 62
$ws1.Cells.Item(112,1).Value = "08:39:38"
$ws1.Cells.Item(112,3).Value = "10_OLMOS"
$ws1.Cells.Item(112,4).Value = 63
$ws1.Cells.Item(113,1).Value = "09:28:24"
$ws1.Cells.Item(113,2).Value = "09:42"
$ws1.Cells.Item(113,3).Value = "215C_EL PATO"
$ws1.Cells.Item(113,4).Value = 14
$ws1.Cells.Item(114,1).Value = "09:28:24"
$ws1.Cells.Item(114,2).Value = "09:43"
$ws1.Cells.Item(114,3).Value = "14_ABASTO"
$ws1.Cells.Item(114,4).Value = 15
$ws1.Cells.Item(115,1).Value = "09:28:24"
$ws1.Cells.Item(115,2).Value = "09:46"
$ws1.Cells.Item(115,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(115,4).Value = 18
$ws1.Cells.Item(116,1).Value = "09:28:24"
$ws1.Cells.Item(116,2).Value = "09:47"
$ws1.Cells.Item(116,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(116,4).Value = 19
$ws1.Cells.Item(117,1).Value = "09:28:24"
$ws1.Cells.Item(117,2).Value = "09:52"
$ws1.Cells.Item(117,4).Value = 24
$ws1.Cells.Item(118,1).Value = "09:28:24"
$ws1.Cells.Item(118,2).Value = "09:53"
$ws1.Cells.Item(118,3).Value = "10_OLMOS"
$ws1.Cells.Item(118,4).Value = 25
$ws1.Cells.Item(119,1).Value = "09:28:24"
$ws1.Cells.Item(119,2).Value = "10:02"
$ws1.Cells.Item(119,3).Value = "17_ROMERO"
$ws1.Cells.Item(119,4).Value = 34
$ws1.Cells.Item(120,1).Value = "09:28:24"
$ws1.Cells.Item(120,2).Value = "10:03"
$ws1.Cells.Item(120,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(120,4).Value = 35
$ws1.Cells.Item(121,1).Value = "09:28:24"
$ws1.Cells.Item(121,2).Value = "10:10"
$ws1.Cells.Item(121,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(121,4).Value = 42
$ws1.Cells.Item(122,1).Value = "09:28:24"
$ws1.Cells.Item(122,2).Value = "10:12"
$ws1.Cells.Item(122,3).Value = "15_ABASTO"
$ws1.Cells.Item(122,4).Value = 44
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "09:28:24"
$ws1.Cells.Item(123,2).Value = "10:13"
$ws1.Cells.Item(123,3).Value = "10_OLMOS"
$ws1.Cells.Item(123,4).Value = 45
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "09:28:24"
$ws1.Cells.Item(124,2).Value = "10:21"
$ws1.Cells.Item(124,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(124,4).Value = 53
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "09:28:24"
$ws1.Cells.Item(125,2).Value = "10:23"
$ws1.Cells.Item(125,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(125,4).Value = 55
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "09:28:24"
$ws1.Cells.Item(126,2).Value = "10:26"
$ws1.Cells.Item(126,3).Value = "215A_EL PATO"
$ws1.Cells.Item(126,4).Value = 58
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "09:28:24"
$ws1.Cells.Item(127,2).Value = "10:34"
$ws1.Cells.Item(127,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(127,4).Value = 66
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "09:28:24"
$ws1.Cells.Item(128,2).Value = "10:42"
$ws1.Cells.Item(128,3).Value = "17_ROMERO"
$ws1.Cells.Item(128,4).Value = 74
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "09:28:24"
$ws1.Cells.Item(129,2).Value = "10:43"
$ws1.Cells.Item(129,3).Value = "14_ABASTO"
$ws1.Cells.Item(129,4).Value = 75
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "09:28:24"
$ws1.Cells.Item(130,2).Value = "10:54"
$ws1.Cells.Item(130,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(130,4).Value = 86
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "09:28:24"
$ws1.Cells.Item(131,2).Value = "11:02"
$ws1.Cells.Item(131,3).Value = "215C_EL PATO"
$ws1.Cells.Item(131,4).Value = 94
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "09:28:24"
$ws1.Cells.Item(132,2).Value = "11:06"
$ws1.Cells.Item(132,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(132,4).Value = 98
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "09:28:24"
$ws1.Cells.Item(133,2).Value = "11:19"
$ws1.Cells.Item(133,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(133,4).Value = 111
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "09:28:24"
$ws1.Cells.Item(134,2).Value = "11:21"
$ws1.Cells.Item(134,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(134,4).Value = 113
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(135,1).Value = "09:28:24"
$ws1.Cells.Item(135,2).Value = "11:27"
$ws1.Cells.Item(135,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(135,4).Value = 119
$ws1.Cells.Item(135,5).Value = "LP1912"

# ---- Sheet 2 ----
$ws2.Cells.Item(2,1).Value = "Última actualización: 09:28:24"
$ws2.Cells.Item(3,1).Value = "Total filas: 23"
$ws2.Cells.Item(26,1).Value = "09:28:24"
$ws2.Cells.Item(26,4).Value = 14
$ws2.Cells.Item(27,1).Value = "09:28:24"
$ws2.Cells.Item(27,4).Value = 58
$ws2.Cells.Item(28,1).Value = "09:28:24"
$ws2.Cells.Item(28,2).Value = "11:02"
$ws2.Cells.Item(28,3).Value = "215C_EL PATO"
$ws2.Cells.Item(28,4).Value = 94
$ws2.Cells.Item(28,5).Value = "LP1912"

# ---- Sheet 3 ----
$ws3.Cells.Item(2,1).Value = "Última actualización: 09:28:24"
$ws3.Cells.Item(3,1).Value = "Total filas: 27"
$ws3.Cells.Item(30,1).Value = "09:28:24"
$ws3.Cells.Item(30,4).Value = 35
$ws3.Cells.Item(31,1).Value = "09:28:24"
$ws3.Cells.Item(31,2).Value = "10:54"
$ws3.Cells.Item(31,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(31,4).Value = 86
$ws3.Cells.Item(31,5).Value = "L6173"
$ws3.Cells.Item(32,1).Value = "09:28:24"
$ws3.Cells.Item(32,2).Value = "11:14"
$ws3.Cells.Item(32,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(32,4).Value = 106
$ws3.Cells.Item(32,5).Value = "L6203"
